$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold/centered/bordered) onto the new header cells I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF) for rows 2-60
$data = @(
    @(2, 6, 6),
    @(3, 7, 7),
    @(4, 7, 7),
    @(5, 9, 9),
    @(6, 6, 6),
    @(7, 8, 8),
    @(8, 8, 8),
    @(9, 4, 5),
    @(10, 6, 6),
    @(11, 7, 7),
    @(12, 7, 7),
    @(13, 7, 7),
    @(14, 7, 7),
    @(15, 7, 7),
    @(16, 6, 6),
    @(17, 7, 7),
    @(18, 6, 6),
    @(19, 8, 8),
    @(20, 6, 6),
    @(21, 7, 7),
    @(22, 7, 7),
    @(23, 7, 7),
    @(24, 7, 7),
    @(25, 6, 7),
    @(26, 8, 8),
    @(27, 7, 7),
    @(28, 7, 8),
    @(29, 7, 7),
    @(30, 7, 7),
    @(31, 7, 8),
    @(32, 8, 9),
    @(33, 7, 8),
    @(34, 7, 7),
    @(35, 7, 7),
    @(36, 7, 7),
    @(37, 6, 7),
    @(38, 9, 9),
    @(39, 7, 7),
    @(40, 7, 7),
    @(41, 8, 8),
    @(42, 6, 6),
    @(43, 8, 8),
    @(44, 6, 6),
    @(45, 8, 8),
    @(46, 4, 5),
    @(47, 7, 7),
    @(48, 8, 8),
    @(49, 7, 7),
    @(50, 6, 7),
    @(51, 8, 9),
    @(52, 7, 7),
    @(53, 9, 9),
    @(54, 6, 7),
    @(55, 8, 8),
    @(56, 9, 9),
    @(57, 8, 8),
    @(58, 5, 5),
    @(59, 5, 5),
    @(60, 5, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
